$d = $word.ActiveDocument

$replacements = @(
    @{old = "326×5=1630"; new = "741×2=1482"},
    @{old = "378×9=3402"; new = "692×8=5536"},
    @{old = "815×3=2445"; new = "754×2=1508"},
    @{old = "477×2=954"; new = "701×6=4206"},
    @{old = "411×2=822"; new = "998×5=4990"},
    @{old = "512×9=4608"; new = "472×9=4248"},
    @{old = "811×2=1622"; new = "847×8=6776"},
    @{old = "226×5=1130"; new = "296×7=2072"},
    @{old = "991×3=2973"; new = "618×4=2472"},
    @{old = "952×4=3808"; new = "175×4=700"},
    @{old = "423×2=846"; new = "853×2=1706"},
    @{old = "653×3=1959"; new = "840×9=7560"},
    @{old = "112×3=336"; new = "371×5=1855"},
    @{old = "987×5=4935"; new = "661×4=2644"},
    @{old = "530×2=1060"; new = "857×4=3428"},
    @{old = "559×9=5031"; new = "922×3=2766"},
    @{old = "490×2=980"; new = "716×2=1432"},
    @{old = "141×8=1128"; new = "783×6=4698"},
    @{old = "386×2=772"; new = "106×8=848"},
    @{old = "688×8=5504"; new = "624×9=5616"},
    @{old = "140×3=420"; new = "456×2=912"},
    @{old = "525×4=2100"; new = "872×9=7848"},
    @{old = "963×2=1926"; new = "475×3=1425"},
    @{old = "530×9=4770"; new = "148×5=740"},
    @{old = "102×4=408"; new = "547×9=4923"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $r.new, 2)
}
